$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to keep a plain-text (string) type even when the
# value looks like a number (e.g. "205.74"), matching the workbook's
# original inline-string cells, then drop back to the default ("Normal")
# cell style so no stray number-format style is left attached to the cell.
function Set-TextValue {
    param($ref, $value)
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "26.708.64"
$ws.Range("E2").Value = "  -1.56%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.544.57"
$ws.Range("E3").Value = "  -1.80%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "205.74"
$ws.Range("E5").Value = "  -0.72%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -1.84%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.04%  "

# Row 8 - Solana
Set-TextValue "D8" "21.34"
$ws.Range("E8").Value = "  -4.15%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -1.66%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -1.23%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -1.96%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "1.767.83"
$ws.Range("E12").Value = "  -1.61%  "

# Row 13 - WrappedEther
Set-TextValue "D13" "1.547.41"
$ws.Range("E13").Value = "  -1.50%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -2.90%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.510"
$ws.Range("E15").Value = "  -1.44%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "26.711.69"
$ws.Range("E16").Value = "  -1.64%  "

# Row 17 - Litecoin
Set-TextValue "D17" "61.12"
$ws.Range("E17").Value = "  -1.82%  "

# Row 18 & 19 - swapped: BitcoinCash <-> ShibaInu
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D18" "0.0₃0689"
$ws.Range("E18").Value = "  +0.70%  "

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D19" "212.26"
$ws.Range("E19").Value = "  -1.04%  "

# Row 20 - Chainlink
Set-TextValue "D20" "7.21"
$ws.Range("E20").Value = "  -2.19%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.06%  "

# Row 22 - Uniswap
Set-TextValue "D22" "4.06"
$ws.Range("E22").Value = "  -1.49%  "

# Row 23 - Avalanche
Set-TextValue "D23" "8.92"
$ws.Range("E23").Value = "  -5.67%  "

# Row 24 - Toncoin
Set-TextValue "D24" "1.98"
$ws.Range("E24").Value = "  -2.20%  "

# Row 25 - Monero
Set-TextValue "D25" "152.64"
$ws.Range("E25").Value = "  +0.16%  "

# Row 26 - Cosmos
Set-TextValue "D26" "6.50"
$ws.Range("E26").Value = "  -3.10%  "

# Row 27 - EthereumClassic
$ws.Range("E27").Value = "  -0.53%  "

# Row 28 - BinanceUSD
$ws.Range("E28").Value = "  -0.04%  "

# Row 29 - Stellar
$ws.Range("E29").Value = "  -1.66%  "

# Row 30 & 31 - swapped: PancakeSwap <-> Hedera
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D30" "0.0459"
$ws.Range("E30").Value = "  -0.87%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D31" "1.10"
$ws.Range("E31").Value = "  -1.67%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +0.26%  "

# Row 33 - Maker
Set-TextValue "D33" "1.335.46"
$ws.Range("E33").Value = "  -4.60%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  -0.49%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -3.29%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -0.70%  "

# Row 37 - TrustWalletToken
Set-TextValue "D37" "0.929"
$ws.Range("E37").Value = "  -1.28%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  -0.30%  "

# Row 39 - ImmutableX
Set-TextValue "D39" "0.522"
$ws.Range("E39").Value = "  +1.36%  "

# Row 40 - FraxShare
Set-TextValue "D40" "5.76"
$ws.Range("E40").Value = "  +6.12%  "

# Row 41 - ARBITRUM
Set-TextValue "D41" "0.798"
$ws.Range("E41").Value = "  -2.05%  "

# Row 42 - WEMIXToken
$ws.Range("E42").Value = "  -1.40%  "

# Row 43 - MXToken
$ws.Range("E43").Value = "  -0.21%  "

# Row 44 & 45 - swapped: RenderToken <-> Aave
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D44" "62.54"
$ws.Range("E44").Value = "  -1.87%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D45" "1.74"
$ws.Range("E45").Value = "  -4.74%  "

# Row 46 - RocketPoolETH
Set-TextValue "D46" "1.681.48"
$ws.Range("E46").Value = "  -1.66%  "

# Row 47 - mCoin
$ws.Range("E47").Value = "  -3.59%  "

# Row 48 - Quant
Set-TextValue "D48" "85.84"
$ws.Range("E48").Value = "  +0.38%  "

# Row 49 - Cronos
Set-TextValue "D49" "0.0505"
$ws.Range("E49").Value = "  +2.44%  "

# Row 50 - BabyDogeCoin
$ws.Range("E50").Value = "  -0.22%  "

# Row 51 - Algorand
Set-TextValue "D51" "0.0952"
$ws.Range("E51").Value = "  -0.14%  "
